$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1090.5454
$ws.Range("I53").Value = 1462.625
$ws.Range("K53").Value = 1462.625
$ws.Range("M53").Value = -825.625
$ws.Range("H76").Value = 3745.8462
$ws.Range("I76").Value = 3812.5
$ws.Range("J76").Value = 3639.2
$ws.Range("K76").Value = 3812.5
$ws.Range("L76").Value = 3639.2
$ws.Range("M76").Value = -3497.5
$ws.Range("N76").Value = -4269.2
$ws.Range("H79").Value = 3745.8462
$ws.Range("I79").Value = 3812.5
$ws.Range("J79").Value = 3639.2
$ws.Range("K79").Value = 3812.5
$ws.Range("L79").Value = 3639.2
$ws.Range("M79").Value = -2720.5
$ws.Range("N79").Value = -5823.2
$ws.Range("H98").Value = 5518.706
$ws.Range("I98").Value = 5676.125
$ws.Range("K98").Value = 5676.125
$ws.Range("M98").Value = -4178.125
$ws.Range("H113").Value = 3138.4
$ws.Range("J113").Value = 3170.6667
$ws.Range("L113").Value = 3170.6667
$ws.Range("N113").Value = -9678.6667
$ws.Range("H118").Value = 1101.8
$ws.Range("J118").Value = 1800
$ws.Range("L118").Value = 5400
$ws.Range("N118").Value = -8714
$ws.Range("H122").Value = 5518.706
$ws.Range("I122").Value = 5676.125
$ws.Range("K122").Value = 17028.375
$ws.Range("M122").Value = -14578.375
$ws.Range("H137").Value = 1408.973
$ws.Range("I137").Value = 1021.88
$ws.Range("J137").Value = 2215.4167
$ws.Range("K137").Value = 3065.64
$ws.Range("L137").Value = 6646.250100000001
$ws.Range("M137").Value = -515.6399999999999
$ws.Range("N137").Value = -11746.2501
$ws.Range("H138").Value = 1792.2716
$ws.Range("J138").Value = 1746.3867
$ws.Range("L138").Value = 5239.1601
$ws.Range("N138").Value = -15519.1601

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3787.8984
$ws.Range("I32").Value = 3639.2031
$ws.Range("K32").Value = 3639.2031
$ws.Range("M32").Value = -3352.2031
$ws.Range("H45").Value = 1300.4667
$ws.Range("I45").Value = 1290.3636
$ws.Range("J45").Value = 1328.25
$ws.Range("K45").Value = 1290.3636
$ws.Range("L45").Value = 1328.25
$ws.Range("M45").Value = -913.3635999999999
$ws.Range("N45").Value = -2082.25
$ws.Range("H61").Value = 66667836
$ws.Range("I61").Value = 83334230
$ws.Range("J61").Value = 2260
$ws.Range("K61").Value = 83334230
$ws.Range("L61").Value = 2260
$ws.Range("M61").Value = -83334018
$ws.Range("N61").Value = -2684
$ws.Range("H74").Value = 2264.125
$ws.Range("I74").Value = 1274.75
$ws.Range("J74").Value = 3253.5
$ws.Range("K74").Value = 1274.75
$ws.Range("L74").Value = 3253.5
$ws.Range("M74").Value = -400.75
$ws.Range("N74").Value = -5001.5
$ws.Range("H77").Value = 2264.125
$ws.Range("I77").Value = 1274.75
$ws.Range("J77").Value = 3253.5
$ws.Range("K77").Value = 6373.75
$ws.Range("L77").Value = 16267.5
$ws.Range("M77").Value = -2005.75
$ws.Range("N77").Value = -25003.5
$ws.Range("H122").Value = 1506.2273
$ws.Range("I122").Value = 1108.3572
$ws.Range("J122").Value = 2202.5
$ws.Range("K122").Value = 3325.0716
$ws.Range("L122").Value = 6607.5
$ws.Range("M122").Value = -875.0715999999998
$ws.Range("N122").Value = -11507.5
$ws.Range("H132").Value = 3814.3333
$ws.Range("I132").Value = 3902.3
$ws.Range("J132").Value = 3704.375
$ws.Range("K132").Value = 11706.9
$ws.Range("L132").Value = 11113.125
$ws.Range("M132").Value = -9176.900000000001
$ws.Range("N132").Value = -16173.125
$ws.Range("H136").Value = 66667836
$ws.Range("I136").Value = 83334230
$ws.Range("J136").Value = 2260
$ws.Range("K136").Value = 250002690
$ws.Range("L136").Value = 6780
$ws.Range("M136").Value = -250000140
$ws.Range("N136").Value = -11880

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 29855
$ws.Range("I57").Value = 27000
$ws.Range("J57").Value = 30262.857
$ws.Range("K57").Value = 27000
$ws.Range("L57").Value = 30262.857
$ws.Range("M57").Value = -26280
$ws.Range("N57").Value = -31702.857
$ws.Range("H107").Value = 1394.1538
$ws.Range("I107").Value = 951.375
$ws.Range("K107").Value = 951.375
$ws.Range("M107").Value = 968.625
$ws.Range("H134").Value = 4055.5293
$ws.Range("I134").Value = 824.55554
$ws.Range("J134").Value = 16517.857
$ws.Range("K134").Value = 2473.66662
$ws.Range("L134").Value = 49553.571
$ws.Range("M134").Value = 61.33338000000003
$ws.Range("N134").Value = -54623.571
$ws.Range("H136").Value = 29855
$ws.Range("I136").Value = 27000
$ws.Range("J136").Value = 30262.857
$ws.Range("K136").Value = 27000
$ws.Range("L136").Value = 30262.857
$ws.Range("M136").Value = -21900
$ws.Range("N136").Value = -40462.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1676.8948
$ws.Range("I58").Value = 1403.9333
$ws.Range("K58").Value = 1403.9333
$ws.Range("M58").Value = -1200.9333
$ws.Range("H74").Value = 28666.666
$ws.Range("I74").Value = 20000
$ws.Range("K74").Value = 20000
$ws.Range("M74").Value = -19126
$ws.Range("H77").Value = 28666.666
$ws.Range("I77").Value = 20000
$ws.Range("K77").Value = 60000
$ws.Range("M77").Value = -55632
$ws.Range("H99").Value = 1322.7273
$ws.Range("I99").Value = 1282
$ws.Range("K99").Value = 1282
$ws.Range("M99").Value = 216
$ws.Range("H114").Value = 23999.5
$ws.Range("J114").Value = 23999.5
$ws.Range("L114").Value = 23999.5
$ws.Range("N114").Value = -32677.5
$ws.Range("H126").Value = 1322.7273
$ws.Range("I126").Value = 1282
$ws.Range("K126").Value = 3846
$ws.Range("M126").Value = -1376
$ws.Range("H131").Value = 14037
$ws.Range("J131").Value = 23778
$ws.Range("L131").Value = 23778
$ws.Range("N131").Value = -33858
$ws.Range("H136").Value = 1676.8948
$ws.Range("I136").Value = 1403.9333
$ws.Range("K136").Value = 4211.7999
$ws.Range("M136").Value = -1661.7999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 637.9167
$ws.Range("I113").Value = 579
$ws.Range("J113").Value = 680
$ws.Range("K113").Value = 1737
$ws.Range("L113").Value = 2040
$ws.Range("M113").Value = 433
$ws.Range("N113").Value = -6380
$ws.Range("H121").Value = 528.5833
$ws.Range("I121").Value = 219
$ws.Range("J121").Value = 962
$ws.Range("K121").Value = 657
$ws.Range("L121").Value = 2886
$ws.Range("M121").Value = 653
$ws.Range("N121").Value = -5506
$ws.Range("H132").Value = 830.7368
$ws.Range("I132").Value = 818
$ws.Range("J132").Value = 858.3333
$ws.Range("K132").Value = 7362
$ws.Range("L132").Value = 7724.9997
$ws.Range("M132").Value = -4832
$ws.Range("N132").Value = -12784.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 58000
$ws.Range("J104").Value = 58000
$ws.Range("L104").Value = 58000
$ws.Range("N104").Value = -64988
$ws.Range("H122").Value = 2403.7222
$ws.Range("I122").Value = 2528.9092
$ws.Range("J122").Value = 2207
$ws.Range("K122").Value = 7586.7276
$ws.Range("L122").Value = 6621
$ws.Range("M122").Value = -5136.7276
$ws.Range("N122").Value = -11521
$ws.Range("H126").Value = 1677.3334
$ws.Range("I126").Value = 1539.25
$ws.Range("J126").Value = 1953.5
$ws.Range("K126").Value = 4617.75
$ws.Range("L126").Value = 5860.5
$ws.Range("M126").Value = -2147.75
$ws.Range("N126").Value = -10800.5
$ws.Range("H132").Value = 4263.2666
$ws.Range("I132").Value = 5787.1665
$ws.Range("K132").Value = 17361.4995
$ws.Range("M132").Value = -14831.4995
$ws.Range("H134").Value = 19316
$ws.Range("J134").Value = 19316
$ws.Range("L134").Value = 57948
$ws.Range("N134").Value = -63018

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1397.8572
$ws.Range("I7").Value = 1450
$ws.Range("J7").Value = 1328.3334
$ws.Range("K7").Value = 1450
$ws.Range("L7").Value = 1328.3334
$ws.Range("M7").Value = -1338
$ws.Range("N7").Value = -1552.3334
$ws.Range("H22").Value = 900
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 1400
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 1400
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -1990
$ws.Range("H27").Value = 900
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 1400
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 1400
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -1614
$ws.Range("H46").Value = 4740.1
$ws.Range("I46").Value = 867
$ws.Range("J46").Value = 6400
$ws.Range("K46").Value = 867
$ws.Range("L46").Value = 6400
$ws.Range("M46").Value = -679
$ws.Range("N46").Value = -6776
$ws.Range("H61").Value = 967.35297
$ws.Range("I61").Value = 957.53845
$ws.Range("K61").Value = 957.53845
$ws.Range("M61").Value = -755.53845
$ws.Range("H100").Value = 1392.3334
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H106").Value = 34200
$ws.Range("J106").Value = 34200
$ws.Range("L106").Value = 34200
$ws.Range("N106").Value = -36724
$ws.Range("H113").Value = 967.35297
$ws.Range("I113").Value = 957.53845
$ws.Range("K113").Value = 957.53845
$ws.Range("M113").Value = 1212.46155
$ws.Range("H122").Value = 22729316
$ws.Range("I122").Value = 31251910
$ws.Range("K122").Value = 93755730
$ws.Range("M122").Value = -93753280
$ws.Range("H126").Value = 1397.8572
$ws.Range("I126").Value = 1450
$ws.Range("J126").Value = 1328.3334
$ws.Range("K126").Value = 4350
$ws.Range("L126").Value = 3985.0002
$ws.Range("M126").Value = -1880
$ws.Range("N126").Value = -8925.0002
$ws.Range("H137").Value = 34580
$ws.Range("J137").Value = 34580
$ws.Range("L137").Value = 34580
$ws.Range("N137").Value = -44780

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1216.7142
$ws.Range("I136").Value = 1162.4546
$ws.Range("K136").Value = 3487.3638
$ws.Range("M136").Value = -937.3638000000001
